$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BASIC")

# Insert a new row above row 22, shifting existing rows (22-38) down to (23-39)
$ws.Rows.Item(22).Insert(-4121)  # -4121 = xlShiftDown

# Copy formatting from the row above (row 20, which holds F20/G20/H20 in the
# same "TOURNAMENT" mini-table, and nothing else) onto the newly inserted
# row's F:H cells only, so the new cells match the surrounding styling
# (fonts / borders) without touching the rest of the (now blank) row.
$ws.Range("F20:H20").Copy()
$ws.Range("F22:H22").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row's content
$ws.Range("F22").Value = "TOURNAMENT_TOTAL_WINNER"
$ws.Range("G22").Value = "NUMBER"
$ws.Range("H22").Value = "COUNT OF TOTAL NUMBER OF PRIZE IN THAT EVENT"

# Update the view so it matches the post-edit selection/scroll position
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("H22").Select()
